# "Label translations" sheet - touch up a handful of the translation/code
# labels ahead of exporting covariance matrices & coefficient means for the
# post-election priors.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C34: fix the short-code label for the "university entrance qualification"
# row (was the old/shorter name).
$ws.Range("C34").Value = "Educ_UnivQualif_Pct"

# B27: replace the German label for GDP per capita with a footnote marker.
$ws.Range("B27").Value = "§"

# C6: fix typo in short code ("Cout" -> "Count").
$ws.Range("C6").Value = "Pop_Total_Count"

# Leave the sheet scrolled to the top with C7 selected, as last saved.
$ws.Range("C7").Select()
